$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2939
$ws.Cells.Item(40, 9).Value = 2939
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 2939
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -2764
$ws.Cells.Item(40, 14).ClearContents()

# ALC row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 502.54544
$ws.Cells.Item(53, 9).Value = 512.25
$ws.Cells.Item(53, 10).Value = 476.66666
$ws.Cells.Item(53, 11).Value = 512.25
$ws.Cells.Item(53, 12).Value = 476.66666
$ws.Cells.Item(53, 13).Value = 124.75
$ws.Cells.Item(53, 14).Value = -1750.66666

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2980.84
$ws.Cells.Item(62, 9).Value = 2500.2942
$ws.Cells.Item(62, 11).Value = 2500.2942
$ws.Cells.Item(62, 13).Value = -1876.2942

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2980.84
$ws.Cells.Item(65, 9).Value = 2500.2942
$ws.Cells.Item(65, 11).Value = 12501.471
$ws.Cells.Item(65, 13).Value = -9381.471

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 134401
$ws.Cells.Item(86, 9).Value = 1601.5
$ws.Cells.Item(86, 10).Value = 400000
$ws.Cells.Item(86, 11).Value = 1601.5
$ws.Cells.Item(86, 12).Value = 400000
$ws.Cells.Item(86, 13).Value = -478.5
$ws.Cells.Item(86, 14).Value = -402246

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 134401
$ws.Cells.Item(89, 9).Value = 1601.5
$ws.Cells.Item(89, 10).Value = 400000
$ws.Cells.Item(89, 11).Value = 8007.5
$ws.Cells.Item(89, 12).Value = 2000000
$ws.Cells.Item(89, 13).Value = -2391.5
$ws.Cells.Item(89, 14).Value = -2011232

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1170.7727
$ws.Cells.Item(112, 9).Value = 850
$ws.Cells.Item(112, 10).Value = 1202.85
$ws.Cells.Item(112, 11).Value = 2550
$ws.Cells.Item(112, 12).Value = 3608.55
$ws.Cells.Item(112, 13).Value = -1442
$ws.Cells.Item(112, 14).Value = -5824.549999999999

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 2500
$ws.Cells.Item(116, 9).Value = 2000
$ws.Cells.Item(116, 10).Value = 3250
$ws.Cells.Item(116, 11).Value = 2000
$ws.Cells.Item(116, 12).Value = 3250
$ws.Cells.Item(116, 13).Value = 1442
$ws.Cells.Item(116, 14).Value = -10134

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1542.8572
$ws.Cells.Item(132, 9).Value = 1353.3846
$ws.Cells.Item(132, 10).Value = 4006
$ws.Cells.Item(132, 11).Value = 4060.1538
$ws.Cells.Item(132, 12).Value = 12018
$ws.Cells.Item(132, 13).Value = -1530.1538
$ws.Cells.Item(132, 14).Value = -17078

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2561.7
$ws.Cells.Item(61, 10).Value = 1400.625
$ws.Cells.Item(61, 12).Value = 1400.625
$ws.Cells.Item(61, 14).Value = -1824.625

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 3433.6667
$ws.Cells.Item(74, 9).Value = 4150.5
$ws.Cells.Item(74, 11).Value = 4150.5
$ws.Cells.Item(74, 13).Value = -3276.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 3433.6667
$ws.Cells.Item(77, 9).Value = 4150.5
$ws.Cells.Item(77, 11).Value = 20752.5
$ws.Cells.Item(77, 13).Value = -16384.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3591.3
$ws.Cells.Item(132, 9).Value = 2600
$ws.Cells.Item(132, 10).Value = 3766.2354
$ws.Cells.Item(132, 11).Value = 7800
$ws.Cells.Item(132, 12).Value = 11298.7062
$ws.Cells.Item(132, 13).Value = -5270
$ws.Cells.Item(132, 14).Value = -16358.7062

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 2561.7
$ws.Cells.Item(136, 10).Value = 1400.625
$ws.Cells.Item(136, 12).Value = 4201.875
$ws.Cells.Item(136, 14).Value = -9301.875

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1998.7894
$ws.Cells.Item(20, 9).Value = 1797.1818
$ws.Cells.Item(20, 10).Value = 2276
$ws.Cells.Item(20, 11).Value = 1797.1818
$ws.Cells.Item(20, 12).Value = 2276
$ws.Cells.Item(20, 13).Value = -1550.1818
$ws.Cells.Item(20, 14).Value = -2770

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2631.92
$ws.Cells.Item(58, 9).Value = 1979.8
$ws.Cells.Item(58, 11).Value = 1979.8
$ws.Cells.Item(58, 13).Value = -1776.8

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 6252.625
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 6252.625
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 18757.875
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -23817.875

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2631.92
$ws.Cells.Item(136, 9).Value = 1979.8
$ws.Cells.Item(136, 11).Value = 5939.4
$ws.Cells.Item(136, 13).Value = -3389.4

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 135.3
$ws.Cells.Item(23, 9).Value = 150
$ws.Cells.Item(23, 10).Value = 129
$ws.Cells.Item(23, 11).Value = 450
$ws.Cells.Item(23, 12).Value = 387
$ws.Cells.Item(23, 13).Value = -215
$ws.Cells.Item(23, 14).Value = -857

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 6660.7646
$ws.Cells.Item(131, 9).Value = 315.55554
$ws.Cells.Item(131, 10).Value = 13799.125
$ws.Cells.Item(131, 11).Value = 946.66662
$ws.Cells.Item(131, 12).Value = 41397.375
$ws.Cells.Item(131, 13).Value = 4093.33338
$ws.Cells.Item(131, 14).Value = -51477.375

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 1909.963
$ws.Cells.Item(136, 9).Value = 1816.0769
$ws.Cells.Item(136, 10).Value = 1997.1428
$ws.Cells.Item(136, 11).Value = 5448.2307
$ws.Cells.Item(136, 12).Value = 5991.428400000001
$ws.Cells.Item(136, 13).Value = -348.2307000000001
$ws.Cells.Item(136, 14).Value = -16191.4284

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(138, 8).Value = 1508.6957
$ws.Cells.Item(138, 9).Value = 973.63635
$ws.Cells.Item(138, 10).Value = 1999.1666
$ws.Cells.Item(138, 11).Value = 2920.90905
$ws.Cells.Item(138, 12).Value = 5997.4998
$ws.Cells.Item(138, 13).Value = 2219.09095
$ws.Cells.Item(138, 14).Value = -16277.4998

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 2177.0588
$ws.Cells.Item(140, 9).Value = 2430
$ws.Cells.Item(140, 10).Value = 2000
$ws.Cells.Item(140, 11).Value = 7290
$ws.Cells.Item(140, 12).Value = 6000
$ws.Cells.Item(140, 13).Value = -2110
$ws.Cells.Item(140, 14).Value = -16360

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 3877.5
$ws.Cells.Item(141, 9).Value = 1530
$ws.Cells.Item(141, 10).Value = 4090.9092
$ws.Cells.Item(141, 11).Value = 4590
$ws.Cells.Item(141, 12).Value = 12272.7276
$ws.Cells.Item(141, 13).Value = 590
$ws.Cells.Item(141, 14).Value = -22632.7276

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2369
$ws.Cells.Item(122, 9).Value = 2142.0322
$ws.Cells.Item(122, 10).Value = 2955.3333
$ws.Cells.Item(122, 11).Value = 6426.096600000001
$ws.Cells.Item(122, 12).Value = 8865.999899999999
$ws.Cells.Item(122, 13).Value = -3976.096600000001
$ws.Cells.Item(122, 14).Value = -13765.9999

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3596.5625
$ws.Cells.Item(132, 9).Value = 3394.4167
$ws.Cells.Item(132, 10).Value = 4203
$ws.Cells.Item(132, 11).Value = 10183.2501
$ws.Cells.Item(132, 12).Value = 12609
$ws.Cells.Item(132, 13).Value = -7653.250100000001
$ws.Cells.Item(132, 14).Value = -17669

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2100.0527
$ws.Cells.Item(132, 9).Value = 2040.3334
$ws.Cells.Item(132, 10).Value = 2246.6365
$ws.Cells.Item(132, 11).Value = 6121.0002
$ws.Cells.Item(132, 12).Value = 6739.9095
$ws.Cells.Item(132, 13).Value = -3591.0002
$ws.Cells.Item(132, 14).Value = -11799.9095

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3523.2856
$ws.Cells.Item(136, 9).Value = 1352.6923
$ws.Cells.Item(136, 10).Value = 7050.5
$ws.Cells.Item(136, 11).Value = 4058.0769
$ws.Cells.Item(136, 12).Value = 21151.5
$ws.Cells.Item(136, 13).Value = -1508.0769
$ws.Cells.Item(136, 14).Value = -26251.5

Write-Host "Edits applied"